# Auto-generated Excel COM-interop script to apply Behemoth_Profits.xlsx diff
# Updates cached values in ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: H6, I6, K6, M6
$ws.Range("H6").Value = 12660.9
$ws.Range("I6").Value = 12660.9
$ws.Range("K6").Value = 37982.7
$ws.Range("M6").Value = -37870.7
# Row 38: H38, I38, J38, K38, L38, M38, N38
$ws.Range("H38").Value = 4521.0527
$ws.Range("I38").Value = 574.25
$ws.Range("J38").Value = 11287
$ws.Range("K38").Value = 1722.75
$ws.Range("L38").Value = 33861
$ws.Range("M38").Value = -1350.75
$ws.Range("N38").Value = -34605
# Row 40: H40, J40, L40, N40
$ws.Range("H40").Value = 3659.92
$ws.Range("J40").Value = 4076.158
$ws.Range("L40").Value = 4076.158
$ws.Range("N40").Value = -4426.157999999999
# Row 41: H41, I41, J41, K41, L41, M41, N41
$ws.Range("H41").Value = 912.44446
$ws.Range("I41").Value = 794.5
$ws.Range("J41").Value = 971.4167
$ws.Range("K41").Value = 794.5
$ws.Range("L41").Value = 971.4167
$ws.Range("M41").Value = -354.5
$ws.Range("N41").Value = -1851.4167
# Row 58: H58, J58, L58, N58
$ws.Range("H58").Value = 12151.125
$ws.Range("J58").Value = 18399.8
$ws.Range("L58").Value = 55199.39999999999
$ws.Range("N58").Value = -55499.39999999999
# Row 62: H62, I62, K62, M62
$ws.Range("H62").Value = 6849.9
$ws.Range("I62").Value = 5374.75
$ws.Range("K62").Value = 5374.75
$ws.Range("M62").Value = -4750.75
# Row 65: H65, I65, K65, M65
$ws.Range("H65").Value = 6849.9
$ws.Range("I65").Value = 5374.75
$ws.Range("K65").Value = 26873.75
$ws.Range("M65").Value = -23753.75
# Row 86: H86, J86, L86, N86
$ws.Range("H86").Value = 8256.6
$ws.Range("J86").Value = 7443.75
$ws.Range("L86").Value = 7443.75
$ws.Range("N86").Value = -9689.75
# Row 87: H87, J87, L87, N87
$ws.Range("H87").Value = 86450
$ws.Range("J87").Value = 86450
$ws.Range("L87").Value = 86450
$ws.Range("N87").Value = -88946
# Row 89: H89, J89, L89, N89
$ws.Range("H89").Value = 8256.6
$ws.Range("J89").Value = 7443.75
$ws.Range("L89").Value = 37218.75
$ws.Range("N89").Value = -48450.75
# Row 90: H90, J90, L90, N90
$ws.Range("H90").Value = 86450
$ws.Range("J90").Value = 86450
$ws.Range("L90").Value = 259350
$ws.Range("N90").Value = -271830
# Row 135: H135, I135, K135, M135
$ws.Range("H135").Value = 8333.375
$ws.Range("I135").Value = 2485.6667
$ws.Range("K135").Value = 22371.0003
$ws.Range("M135").Value = -19836.0003
# Row 137: H137, I137, K137, M137
$ws.Range("H137").Value = 16996.822
$ws.Range("I137").Value = 39373.633
$ws.Range("K137").Value = 118120.899
$ws.Range("M137").Value = -115570.899
# Row 138: H138, J138, L138, N138
$ws.Range("H138").Value = 4518.875
$ws.Range("J138").Value = 5187.804
$ws.Range("L138").Value = 15563.412
$ws.Range("N138").Value = -25843.412

$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32, I32, J32, K32, L32, M32, N32
$ws.Range("H32").Value = 7701876.5
$ws.Range("I32").Value = 8071257
$ws.Range("J32").Value = 68009.336
$ws.Range("K32").Value = 8071257
$ws.Range("L32").Value = 68009.336
$ws.Range("M32").Value = -8070970
$ws.Range("N32").Value = -68583.336
# Row 102: H102, I102, K102, M102
$ws.Range("H102").Value = 22511.5
$ws.Range("I102").Value = 28592.334
$ws.Range("K102").Value = 28592.334
$ws.Range("M102").Value = -26970.334
# Row 103: H103, I103, J103, K103, L103, M103, N103
$ws.Range("H103").Value = 52996.668
$ws.Range("I103").Value = 29000
$ws.Range("J103").Value = 64995
$ws.Range("K103").Value = 29000
$ws.Range("L103").Value = 64995
$ws.Range("M103").Value = -27828
$ws.Range("N103").Value = -67339
# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 4908.3335
$ws.Range("I132").Value = 1685.7693
$ws.Range("K132").Value = 5057.3079
$ws.Range("M132").Value = -2527.3079
# Row 135: H135, J135, L135, N135
$ws.Range("H135").Value = 43979
$ws.Range("J135").Value = 43979
$ws.Range("L135").Value = 43979
$ws.Range("N135").Value = -54119

$ws = $wb.Worksheets.Item("BSM")
# Row 100: H100, J100, L100, N100
$ws.Range("H100").Value = 37718.5
$ws.Range("J100").Value = 37718.5
$ws.Range("L100").Value = 37718.5
$ws.Range("N100").Value = -39882.5

$ws = $wb.Worksheets.Item("CRP")
# Row 28: H28, J28, L28, N28
$ws.Range("H28").Value = 36128.777
$ws.Range("J28").Value = 36128.777
$ws.Range("L28").Value = 36128.777
$ws.Range("N28").Value = -36618.777
# Row 31: H31, J31, L31, N31
$ws.Range("H31").Value = 462393.22
$ws.Range("J31").Value = 835236.7
$ws.Range("L31").Value = 835236.7
$ws.Range("N31").Value = -835826.7
# Row 34: H34, J34, L34, N34
$ws.Range("H34").Value = 462393.22
$ws.Range("J34").Value = 835236.7
$ws.Range("L34").Value = 835236.7
$ws.Range("N34").Value = -835640.7
# Row 86: H86, I86, K86, M86
$ws.Range("H86").Value = 4285
$ws.Range("I86").Value = 3999.6667
$ws.Range("K86").Value = 3999.6667
$ws.Range("M86").Value = -2876.6667
# Row 89: H89, I89, K89, M89
$ws.Range("H89").Value = 4285
$ws.Range("I89").Value = 3999.6667
$ws.Range("K89").Value = 19998.3335
$ws.Range("M89").Value = -14382.3335
# Row 109: H109, J109, L109, N109
$ws.Range("H109").Value = 33666.668
$ws.Range("J109").Value = 33666.668
$ws.Range("L109").Value = 33666.668
$ws.Range("N109").Value = -35746.668
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 2962.4707
$ws.Range("I132").Value = 2797.077
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 8391.231
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -5861.231
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("CUL")
# Row 140: H140, I140, K140, M140
$ws.Range("H140").Value = 177503
$ws.Range("I140").Value = 177503
$ws.Range("K140").Value = 532509
$ws.Range("M140").Value = -527329

$ws = $wb.Worksheets.Item("GSM")
# Row 70: H70, I70, J70, K70, L70, M70, N70
$ws.Range("H70").Value = 6347.8
$ws.Range("I70").Value = 5536.5
$ws.Range("J70").Value = 7564.75
$ws.Range("K70").Value = 5536.5
$ws.Range("L70").Value = 7564.75
$ws.Range("M70").Value = -5266.5
$ws.Range("N70").Value = -8104.75
# Row 73: H73, I73, J73, K73, L73, M73, N73
$ws.Range("H73").Value = 6347.8
$ws.Range("I73").Value = 5536.5
$ws.Range("J73").Value = 7564.75
$ws.Range("K73").Value = 5536.5
$ws.Range("L73").Value = 7564.75
$ws.Range("M73").Value = -4600.5
$ws.Range("N73").Value = -9436.75
# Row 96: H96, J96, L96, N96
$ws.Range("H96").Value = 57499.25
$ws.Range("J96").Value = 69999.336
$ws.Range("L96").Value = 69999.336
$ws.Range("N96").Value = -75491.336
# Row 97: H97, I97, J97, K97, L97, M97, N97
$ws.Range("H97").Value = 1926.174
$ws.Range("I97").Value = 1940.15
$ws.Range("J97").Value = 1833
$ws.Range("K97").Value = 1940.15
$ws.Range("L97").Value = 1833
$ws.Range("M97").Value = -1444.15
$ws.Range("N97").Value = -2825
# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value = 11015.857
$ws.Range("I122").Value = 7442.6
$ws.Range("J122").Value = 19949
$ws.Range("K122").Value = 22327.8
$ws.Range("L122").Value = 59847
$ws.Range("M122").Value = -19877.8
$ws.Range("N122").Value = -64747
# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 166670770
$ws.Range("I132").Value = 166670770
$ws.Range("K132").Value = 500012310
$ws.Range("M132").Value = -500009780
# Row 134: H134, J134, L134, N134
$ws.Range("H134").Value = 90163
$ws.Range("J134").Value = 90163
$ws.Range("L134").Value = 270489
$ws.Range("N134").Value = -275559
# Row 141: H141, J141, L141, N141
$ws.Range("H141").Value = 157249.75
$ws.Range("J141").Value = 157249.75
$ws.Range("L141").Value = 157249.75
$ws.Range("N141").Value = -167609.75

$ws = $wb.Worksheets.Item("LTW")
# Row 122: H122, I122, K122, M122
$ws.Range("H122").Value = 6097.926
$ws.Range("I122").Value = 5112.579
$ws.Range("K122").Value = 15337.737
$ws.Range("M122").Value = -12887.737
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 41966.137
$ws.Range("I132").Value = 64573.832
$ws.Range("J132").Value = 4971.727
$ws.Range("K132").Value = 193721.496
$ws.Range("L132").Value = 14915.181
$ws.Range("M132").Value = -191191.496
$ws.Range("N132").Value = -19975.181
# Row 136: H136, I136, J136, K136, L136, M136, N136
$ws.Range("H136").Value = 106334.09
$ws.Range("I136").Value = 148634.42
$ws.Range("J136").Value = 86593.92999999999
$ws.Range("K136").Value = 445903.26
$ws.Range("L136").Value = 259781.79
$ws.Range("M136").Value = -443353.26
$ws.Range("N136").Value = -264881.79

$ws = $wb.Worksheets.Item("WVR")
# Row 14: H14, I14, K14, M14
$ws.Range("H14").Value = 3832.3333
$ws.Range("I14").Value = 2748.5
$ws.Range("K14").Value = 2748.5
$ws.Range("M14").Value = -2580.5
# Row 33: H33, J33, L33, N33
$ws.Range("H33").Value = 28247.5
$ws.Range("J33").Value = 28247.5
$ws.Range("L33").Value = 28247.5
$ws.Range("N33").Value = -28747.5
# Row 36: H36, J36, L36, N36
$ws.Range("H36").Value = 28247.5
$ws.Range("J36").Value = 28247.5
$ws.Range("L36").Value = 28247.5
$ws.Range("N36").Value = -28747.5
# Row 63: H63, J63, L63, N63
$ws.Range("H63").Value = 19373.25
$ws.Range("J63").Value = 18497.666
$ws.Range("L63").Value = 18497.666
$ws.Range("N63").Value = -19745.666
# Row 66: H66, J66, L66, N66
$ws.Range("H66").Value = 19373.25
$ws.Range("J66").Value = 18497.666
$ws.Range("L66").Value = 55492.99800000001
$ws.Range("N66").Value = -61732.99800000001
# Row 98: H98, J98, L98, N98
$ws.Range("H98").Value = 38295
$ws.Range("J98").Value = 38295
$ws.Range("L98").Value = 38295
$ws.Range("N98").Value = -44285
# Row 101: H101, J101, L101, N101
$ws.Range("H101").Value = 43440.6
$ws.Range("J101").Value = 43440.6
$ws.Range("L101").Value = 43440.6
$ws.Range("N101").Value = -49930.6
# Row 122: H122, I122, K122, M122
$ws.Range("H122").Value = 7678.3477
$ws.Range("I122").Value = 3422.9443
$ws.Range("K122").Value = 10268.8329
$ws.Range("M122").Value = -7818.832900000001
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 1002989.4
$ws.Range("I132").Value = 3328.1765
$ws.Range("J132").Value = 6667736.5
$ws.Range("K132").Value = 9984.529500000001
$ws.Range("L132").Value = 20003209.5
$ws.Range("M132").Value = -7454.529500000001
$ws.Range("N132").Value = -20008269.5
# Row 136: H136, I136, K136, M136
$ws.Range("H136").Value = 3955.7334
$ws.Range("I136").Value = 3025.4614
$ws.Range("K136").Value = 9076.3842
$ws.Range("M136").Value = -6526.3842
# Row 139: H139, J139, L139, N139
$ws.Range("H139").Value = 150142.8
$ws.Range("J139").Value = 150142.8
$ws.Range("L139").Value = 150142.8
$ws.Range("N139").Value = -160422.8
